$wb = $excel.ActiveWorkbook

# --- Data sheet: append two new weekly observations (FRED pull refresh) ---
$dataWs = $wb.Worksheets.Item("Data")

# Seed rows 98 and 99 by copying the formatting of the last existing data
# row (97) down, so the new date cells keep the same date number format /
# font / border / alignment as the rest of column A.
$dataWs.Range("A97:B97").Copy($dataWs.Range("A98:B98"))
$dataWs.Range("A97:B97").Copy($dataWs.Range("A99:B99"))

$dataWs.Cells.Item(98, 1).Value = 45147
$dataWs.Cells.Item(98, 2).Value = 3222.856

$dataWs.Cells.Item(99, 1).Value = 45154
$dataWs.Cells.Item(99, 2).Value = 3245.971

# --- SeriesInfo sheet: update metadata to reflect the refreshed pull ---
$infoWs = $wb.Worksheets.Item("SeriesInfo")

# These cells hold plain text (e.g. "2023-08-22"), not real dates, in the
# source workbook. Force a text number format before assigning so Excel
# doesn't auto-convert the date-shaped strings into date serial numbers,
# then clear the cell-level formatting again so the cell keeps the same
# (default) style it had before, matching the source file.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $infoWs.Range("B3") "2023-08-22"
Set-TextValue $infoWs.Range("B4") "2023-08-22"
Set-TextValue $infoWs.Range("B7") "2023-08-16"
Set-TextValue $infoWs.Range("B14") "2023-08-17 15:35:51-05"

$infoWs.Range("B15").Value = 75
